# Auto-generated edit script applying the Ultros_Profits.xlsx cell-value diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H118").Value = 1458.625
$ws.Range("J118").Value = 999.75
$ws.Range("L118").Value = 2999.25
$ws.Range("N118").Value = -6313.25
$ws.Range("H132").Value = 25170.36
$ws.Range("I132").Value = 3963.2942
$ws.Range("K132").Value = 11889.8826
$ws.Range("M132").Value = -9359.882599999999
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88
$ws.Range("H32").Value = 20005418
$ws.Range("I32").Value = 20413650
$ws.Range("K32").Value = 20413650
$ws.Range("M32").Value = -20413363
$ws.Range("H45").Value = 2808
$ws.Range("I45").Value = 2212.25
$ws.Range("K45").Value = 2212.25
$ws.Range("M45").Value = -1835.25
$ws.Range("H74").Value = 2808.1667
$ws.Range("I74").Value = 2534.1875
$ws.Range("K74").Value = 2534.1875
$ws.Range("M74").Value = -1660.1875
$ws.Range("H77").Value = 2808.1667
$ws.Range("I77").Value = 2534.1875
$ws.Range("K77").Value = 12670.9375
$ws.Range("M77").Value = -8302.9375
$ws.Range("H110").Value = 4258.5
$ws.Range("I110").Value = 4258.5
$ws.Range("K110").Value = 4258.5
$ws.Range("M110").Value = -2213.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85
$ws.Range("H64").Value = 925.3
$ws.Range("J64").Value = 985.4286
$ws.Range("L64").Value = 985.4286
$ws.Range("N64").Value = -1435.4286
$ws.Range("H67").Value = 925.3
$ws.Range("J67").Value = 985.4286
$ws.Range("L67").Value = 985.4286
$ws.Range("N67").Value = -2545.4286
$ws.Range("H86").Value = 2368.2917
$ws.Range("I86").Value = 2036.5625
$ws.Range("J86").Value = 3031.75
$ws.Range("K86").Value = 2036.5625
$ws.Range("L86").Value = 3031.75
$ws.Range("M86").Value = -913.5625
$ws.Range("N86").Value = -5277.75
$ws.Range("H89").Value = 2368.2917
$ws.Range("I89").Value = 2036.5625
$ws.Range("J89").Value = 3031.75
$ws.Range("K89").Value = 10182.8125
$ws.Range("L89").Value = 15158.75
$ws.Range("M89").Value = -4566.8125
$ws.Range("N89").Value = -26390.75
$ws.Range("H94").Value = 1730.3478
$ws.Range("I94").Value = 1970.9286
$ws.Range("K94").Value = 1970.9286
$ws.Range("M94").Value = -1519.9286
$ws.Range("H134").Value = 1680.5758
$ws.Range("I134").Value = 1514.3438
$ws.Range("K134").Value = 4543.0314
$ws.Range("M134").Value = -2008.0314
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 789
$ws.Range("I22").Value = 789
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 789
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -439
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 1690.0817
$ws.Range("I31").Value = 1578.1282
$ws.Range("J31").Value = 2126.7
$ws.Range("K31").Value = 1578.1282
$ws.Range("L31").Value = 2126.7
$ws.Range("M31").Value = -1283.1282
$ws.Range("N31").Value = -2716.7
$ws.Range("H34").Value = 1690.0817
$ws.Range("I34").Value = 1578.1282
$ws.Range("J34").Value = 2126.7
$ws.Range("K34").Value = 1578.1282
$ws.Range("L34").Value = 2126.7
$ws.Range("M34").Value = -1376.1282
$ws.Range("N34").Value = -2530.7
$ws.Range("H94").Value = 10626.637
$ws.Range("I94").Value = 20649.4
$ws.Range("K94").Value = 20649.4
$ws.Range("M94").Value = -20198.4
$ws.Range("H107").Value = 11760.158
$ws.Range("I107").Value = 1241.1428
$ws.Range("J107").Value = 17896.25
$ws.Range("K107").Value = 1241.1428
$ws.Range("L107").Value = 17896.25
$ws.Range("M107").Value = 678.8571999999999
$ws.Range("N107").Value = -21736.25
$ws.Range("H132").Value = 2081.2354
$ws.Range("I132").Value = 2047.9375
$ws.Range("K132").Value = 6143.8125
$ws.Range("M132").Value = -3613.8125
$ws.Range("H134").Value = 2849.28
$ws.Range("I134").Value = 2263.25
$ws.Range("K134").Value = 6789.75
$ws.Range("M134").Value = -4254.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 475.85715
$ws.Range("J52").Value = 475.85715
$ws.Range("L52").Value = 1427.57145
$ws.Range("N52").Value = -1959.57145
$ws.Range("H60").Value = 2080
$ws.Range("J60").Value = 3116.6667
$ws.Range("L60").Value = 9350.000100000001
$ws.Range("N60").Value = -9852.000100000001
$ws.Range("H109").Value = 369.33334
$ws.Range("I109").Value = 369.33334
$ws.Range("K109").Value = 1108.00002
$ws.Range("M109").Value = -68.00001999999995
$ws.Range("H113").Value = 1310.4706
$ws.Range("J113").Value = 1301.7
$ws.Range("L113").Value = 3905.1
$ws.Range("N113").Value = -8245.1
$ws.Range("H119").Value = 3812.375
$ws.Range("J119").Value = 6000
$ws.Range("L119").Value = 18000
$ws.Range("N119").Value = -27676
$ws.Range("H122").Value = 707.75
$ws.Range("I122").Value = 799.8
$ws.Range("J122").Value = 642
$ws.Range("K122").Value = 7198.2
$ws.Range("L122").Value = 5778
$ws.Range("M122").Value = -4748.2
$ws.Range("N122").Value = -10678
$ws.Range("H131").Value = 2869.423
$ws.Range("J131").Value = 4129.4165
$ws.Range("L131").Value = 12388.2495
$ws.Range("N131").Value = -22468.2495
$ws.Range("H132").Value = 2664.8333
$ws.Range("I132").Value = 2499.75
$ws.Range("K132").Value = 22497.75
$ws.Range("M132").Value = -19967.75
$ws.Range("H138").Value = 5465.731
$ws.Range("I138").Value = 3185.75
$ws.Range("K138").Value = 9557.25
$ws.Range("M138").Value = -4417.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 13348.75
$ws.Range("J20").Value = 13348.75
$ws.Range("L20").Value = 13348.75
$ws.Range("N20").Value = -13838.75
$ws.Range("H70").Value = 529402.75
$ws.Range("I70").Value = 529402.75
$ws.Range("K70").Value = 529402.75
$ws.Range("M70").Value = -529132.75
$ws.Range("H73").Value = 529402.75
$ws.Range("I73").Value = 529402.75
$ws.Range("K73").Value = 529402.75
$ws.Range("M73").Value = -528466.75
$ws.Range("H97").Value = 313.45834
$ws.Range("I97").Value = 277.3684
$ws.Range("K97").Value = 277.3684
$ws.Range("M97").Value = 218.6316
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("H122").Value = 5252.3125
$ws.Range("I122").Value = 4185.364
$ws.Range("K122").Value = 12556.092
$ws.Range("M122").Value = -10106.092
$ws.Range("H132").Value = 4483
$ws.Range("I132").Value = 4332.6665
$ws.Range("J132").Value = 4633.3335
$ws.Range("K132").Value = 12997.9995
$ws.Range("L132").Value = 13900.0005
$ws.Range("M132").Value = -10467.9995
$ws.Range("N132").Value = -18960.0005
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1500
$ws.Range("I9").Value = 1000
$ws.Range("K9").Value = 1000
$ws.Range("M9").Value = -776
$ws.Range("H68").Value = 2935.7646
$ws.Range("I68").Value = 2969.0908
$ws.Range("J68").Value = 2874.6667
$ws.Range("K68").Value = 2969.0908
$ws.Range("L68").Value = 2874.6667
$ws.Range("M68").Value = -2220.0908
$ws.Range("N68").Value = -4372.6667
$ws.Range("H71").Value = 2935.7646
$ws.Range("I71").Value = 2969.0908
$ws.Range("J71").Value = 2874.6667
$ws.Range("K71").Value = 14845.454
$ws.Range("L71").Value = 14373.3335
$ws.Range("M71").Value = -11101.454
$ws.Range("N71").Value = -21861.3335
$ws.Range("H93").Value = 8208.166999999999
$ws.Range("I93").Value = 8007.636
$ws.Range("K93").Value = 8007.636
$ws.Range("M93").Value = -6759.636
$ws.Range("H132").Value = 3305.9614
$ws.Range("I132").Value = 2390.1052
$ws.Range("K132").Value = 7170.3156
$ws.Range("M132").Value = -4640.3156
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2500
$ws.Range("K62").Value = 2500
$ws.Range("M62").Value = -1876
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2500
$ws.Range("K65").Value = 12500
$ws.Range("M65").Value = -9380

Write-Output "Applied 210 cell edits"
